$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.457073490363895
$ws.Range("B2").Value = -1.247247429609309

$ws.Range("A3").Value = -0.5353822030659443
$ws.Range("B3").Value = -0.6694704817216622

$ws.Range("A4").Value = -0.8498925758554317
$ws.Range("B4").Value = -0.6843678232413823

$ws.Range("A5").Value = -0.7319290046682726
$ws.Range("B5").Value = -0.6349686399606124

$ws.Range("A6").Value = 0.8083947436854791
$ws.Range("B6").Value = 0.6068203705027657

$ws.Range("A7").Value = -0.07452177043058129
$ws.Range("B7").Value = 0.0504749455727569

$ws.Range("A8").Value = 0.7880402694620865
$ws.Range("B8").Value = 0.5459839380625998

$ws.Range("A9").Value = 0.3186855317368091
$ws.Range("B9").Value = 0.2961332645006233
